$d = $word.ActiveDocument

# Locate the paragraph that still ends with the requirements text, then
# remove the three paragraphs that followed it in the old footer block:
#   - a blank paragraph
#   - "Ver no Jupiter Salvar em pdf Salvar em docx"
#   - the "(c) 2020 ... Powered by Jekyll ..." copyright line
# This collapses them away while leaving the paragraph mark of the
# "LOQ4205..." paragraph and the subsequent (already blank / page-break)
# paragraphs untouched.

$marker = "LOQ4205: Sistemas Produtivos II (Requisito fraco)"
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq $marker) {
        $target = $i
        break
    }
}

if ($target -eq $null) {
    throw "Could not find paragraph with text: $marker"
}

$startPara = $d.Paragraphs.Item($target)

# Find the copyright paragraph (3 paragraphs after the marker) by text,
# walking forward robustly instead of assuming a fixed offset.
$copyrightIdx = $null
for ($j = $target + 1; $j -le $d.Paragraphs.Count; $j++) {
    $txt = $d.Paragraphs.Item($j).Range.Text
    if ($txt -like "*Powered by Jekyll*") {
        $copyrightIdx = $j
        break
    }
}

if ($copyrightIdx -eq $null) {
    throw "Could not find the copyright paragraph"
}

$endPara = $d.Paragraphs.Item($copyrightIdx)

# Delete from right after the marker paragraph's mark through the end of
# (and including) the copyright paragraph's mark - this removes the blank
# paragraph, the "Ver no Jupiter..." paragraph and the copyright paragraph
# in one shot, including all of their paragraph marks.
$rng = $d.Range($startPara.Range.End, $endPara.Range.End)
$rng.Delete()
